# Add a new "2022-Q1" sheet (fund holdings detail) right before the "总计"
# (total) summary sheet, and prepend a corresponding "2022-Q1" row to the
# "总计" sheet's summary table.

$wb = $excel.ActiveWorkbook

$totalSheetBefore = $wb.Worksheets.Item("总计")
$q4Sheet          = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet right before "总计"
# ---------------------------------------------------------------------
# NOTE: this COM-interop engine resolves worksheet references by
# position, so the handle used for "Before" (as well as any other
# handle previously obtained) is NOT stable across an Add() call - it
# ends up pointing at whatever sheet now occupies that slot (i.e. the
# newly inserted one). Worksheet references must therefore be re-
# fetched (by name) after any Worksheets.Add() call.
$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# Re-fetch a fresh, valid reference to the "总计" sheet now that the
# insertion above has shifted sheet positions around.
$totalSheet = $wb.Worksheets.Item("总计")

# Replicate the header-row style (bold/border) and index-column style
# used by the other quarterly sheets (e.g. "2021-Q4"), so the new sheet
# matches the established look & feel.
$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$q4Sheet.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B-G hold numeric-looking values but are stored as text in this
# workbook's convention, so force a text format before assigning them.
$newSheet.Range("B2:G5").NumberFormat = "@"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "011243"
$newSheet.Range("C2").Value = "万家惠裕回报6个月持有期混合型证券投资基金A"
$newSheet.Range("D2").Value = "4.93"
$newSheet.Range("E2").Value = "23.04"
$newSheet.Range("F2").Value = "0.96"
$newSheet.Range("G2").Value = "0.0473"
$newSheet.Range("H2").Value = 5

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "519656"
$newSheet.Range("C3").Value = "银河灵活配置混合 - A"
$newSheet.Range("D3").Value = "0.72"
$newSheet.Range("E3").Value = "59.27"
$newSheet.Range("F3").Value = "2.99"
$newSheet.Range("G3").Value = "0.0215"
$newSheet.Range("H3").Value = 9

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "519657"
$newSheet.Range("C4").Value = "银河灵活配置混合 - C"
$newSheet.Range("D4").Value = "0.33"
$newSheet.Range("E4").Value = "59.27"
$newSheet.Range("F4").Value = "2.99"
$newSheet.Range("G4").Value = "0.0099"
$newSheet.Range("H4").Value = 9

# Row 5
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "011244"
$newSheet.Range("C5").Value = "万家惠裕回报6个月持有期混合型证券投资基金C"
$newSheet.Range("D5").Value = "0.14"
$newSheet.Range("E5").Value = "23.04"
$newSheet.Range("F5").Value = "0.96"
$newSheet.Range("G5").Value = "0.0013"
$newSheet.Range("H5").Value = 5

# The text-format markup applied above (to keep B-G as text) is only
# needed to influence how the literal values are parsed; drop it again
# afterwards so the cells end up unstyled, just like their counterparts
# on the other quarterly sheets.
$newSheet.Range("B2:G5").ClearFormats()

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q1" summary row to the "总计" sheet
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Copy the index-column style (bold/border) down into the newly
# inserted row.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)    # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.08

# Renumber the index column so it stays a contiguous 0-based sequence.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
